$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 47, shifting existing rows 47-147 down to 48-148.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new data record.
$ws.Range("A47").Value = 10
$ws.Range("B47").Value = "Vega Modelo de Temuco"
$ws.Range("C47").Value = "La Araucanía"
$ws.Range("D47").Value = (Get-Date -Year 2022 -Month 8 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E47").Value = 9
$ws.Range("F47").Value = 100112031
$ws.Range("G47").Value = "Poroto verde"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 40
$ws.Range("K47").Value = 30000
$ws.Range("L47").Value = 30000
$ws.Range("M47").Value = 30000
$ws.Range("N47").Value = "`$/malla 25 kilos"
$ws.Range("O47").Value = "Provincia de Limarí"
$ws.Range("P47").Value = 1200
$ws.Range("Q47").Value = 25
$ws.Range("R47").Value = "Hortaliza"

Write-Output "done"
